$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '32.816.42'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +9.60%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.759.19'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +5.96%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.998'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.02%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '227.01'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +4.42%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.544'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +4.28%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.999'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.06%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '31.73'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +9.81%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '45.25'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +3.24%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.278'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +5.42%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0667'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +8.42%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0918'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.67%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '2.008.11'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +5.77%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.750.48'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +5.28%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.631'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +3.45%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '10.44'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +2.85%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '4.30'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +8.44%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '32.763.58'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +9.32%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '68.74'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +5.69%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '258.17'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +6.16%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.0₃0743'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +4.57%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.996'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.17%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.46'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +4.29%  '

$ws.Range("E24").Value = '  +3.79%  '

$ws.Range("E25").Value = '  -0.81%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '159.58'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.50%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '16.51'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +4.39%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.116'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +4.38%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '6.97'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +3.59%  '

$ws.Range("E30").Value = '  +0.22%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.87'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +13.35%  '

$ws.Range("E32").Value = '  +3.39%  '

$ws.Range("E33").Value = '  +5.27%  '

$ws.Range("E34").Value = '  +8.05%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.557.29'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +7.61%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.79'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +4.31%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.05'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +2.00%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.632'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +9.77%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '84.41'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +5.92%  '

$ws.Range("E40").Value = '  +5.41%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.74'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +3.71%  '

$ws.Range("E42").Value = '  +0.80%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.875'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +2.90%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.08'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +6.47%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0513'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +2.87%  '

$ws.Range("E46").Value = '  +4.43%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '53.97'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +6.74%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.908.54'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +5.77%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.00'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.15%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '5.68'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +5.75%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '95.69'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.47%  '
